# Auto-generated script applying scheduled-runner price/profit updates
# to the Famfrit_Profits leve tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 24999.5
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H69").Value = 6000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null
$ws.Range("H72").Value = 6000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null
$ws.Range("H96").Value = 1794.3889
$ws.Range("I96").Value = 1411.5385
$ws.Range("J96").Value = 2789.8
$ws.Range("K96").Value = 4234.6155
$ws.Range("L96").Value = 8369.400000000001
$ws.Range("M96").Value = -2861.6155
$ws.Range("N96").Value = -11115.4
$ws.Range("H135").Value = 1558.1666
$ws.Range("I135").Value = 1599.75
$ws.Range("K135").Value = 14397.75
$ws.Range("M135").Value = -11862.75
$ws.Range("H137").Value = 14422.322
$ws.Range("I137").Value = 19324.055
$ws.Range("J137").Value = 7635.3076
$ws.Range("K137").Value = 57972.165
$ws.Range("L137").Value = 22905.9228
$ws.Range("M137").Value = -55422.165
$ws.Range("N137").Value = -28005.9228
$ws.Range("H138").Value = 2613
$ws.Range("I138").Value = 1204.2307
$ws.Range("J138").Value = 3576.8948
$ws.Range("K138").Value = 3612.6921
$ws.Range("L138").Value = 10730.6844
$ws.Range("M138").Value = 1527.3079
$ws.Range("N138").Value = -21010.6844
$ws.Range("H141").Value = 10367.5
$ws.Range("I141").Value = 13998.75
$ws.Range("K141").Value = 41996.25
$ws.Range("M141").Value = -36816.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 617.375
$ws.Range("I5").Value = 691.2857
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 691.2857
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -579.2857
$ws.Range("N5").Value = -324
$ws.Range("H74").Value = 31460.54
$ws.Range("I74").Value = 32182.816
$ws.Range("K74").Value = 32182.816
$ws.Range("M74").Value = -31308.816
$ws.Range("H77").Value = 31460.54
$ws.Range("I77").Value = 32182.816
$ws.Range("K77").Value = 160914.08
$ws.Range("M77").Value = -156546.08
$ws.Range("H107").Value = 40000
$ws.Range("I107").Value = 40000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 40000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -36160
$ws.Range("N107").Value = $null
$ws.Range("H122").Value = 3121
$ws.Range("I122").Value = 2725
$ws.Range("K122").Value = 8175
$ws.Range("M122").Value = -5725
$ws.Range("H132").Value = 66620.266
$ws.Range("I132").Value = 4085.3
$ws.Range("J132").Value = 301126.38
$ws.Range("K132").Value = 12255.9
$ws.Range("L132").Value = 903379.14
$ws.Range("M132").Value = -9725.900000000001
$ws.Range("N132").Value = -908439.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 617.375
$ws.Range("I4").Value = 691.2857
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 691.2857
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -576.2857
$ws.Range("N4").Value = -330
$ws.Range("H75").Value = 31664.666
$ws.Range("I75").Value = 7497
$ws.Range("J75").Value = 80000
$ws.Range("K75").Value = 7497
$ws.Range("L75").Value = 80000
$ws.Range("M75").Value = -6561
$ws.Range("N75").Value = -81872
$ws.Range("H78").Value = 31664.666
$ws.Range("I78").Value = 7497
$ws.Range("J78").Value = 80000
$ws.Range("K78").Value = 22491
$ws.Range("L78").Value = 240000
$ws.Range("M78").Value = -17811
$ws.Range("N78").Value = -249360
$ws.Range("H123").Value = 84997.5
$ws.Range("J123").Value = 84997.5
$ws.Range("L123").Value = 84997.5
$ws.Range("N123").Value = -94797.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19068.857
$ws.Range("J41").Value = 31132
$ws.Range("L41").Value = 31132
$ws.Range("N41").Value = -31988
$ws.Range("H50").Value = 14000
$ws.Range("I50").Value = 14000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 14000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -13375
$ws.Range("N50").Value = $null
$ws.Range("H51").Value = 19900
$ws.Range("I51").Value = 19900
$ws.Range("K51").Value = 19900
$ws.Range("M51").Value = -19164
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = $null
$ws.Range("H61").Value = 19900
$ws.Range("I61").Value = 19900
$ws.Range("K61").Value = 19900
$ws.Range("M61").Value = -19552
$ws.Range("H134").Value = 4390.25
$ws.Range("I134").Value = 3853.8333
$ws.Range("K134").Value = 11561.4999
$ws.Range("M134").Value = -9026.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1452.0769
$ws.Range("I5").Value = 175.22223
$ws.Range("J5").Value = 4325
$ws.Range("K5").Value = 525.66669
$ws.Range("L5").Value = 12975
$ws.Range("M5").Value = -413.66669
$ws.Range("N5").Value = -13199
$ws.Range("H21").Value = 541.2857
$ws.Range("I21").Value = 246.8
$ws.Range("J21").Value = 1277.5
$ws.Range("K21").Value = 740.4000000000001
$ws.Range("L21").Value = 3832.5
$ws.Range("M21").Value = -567.4000000000001
$ws.Range("N21").Value = -4178.5
$ws.Range("H122").Value = 1824.9286
$ws.Range("I122").Value = 751
$ws.Range("J122").Value = 2254.5
$ws.Range("K122").Value = 6759
$ws.Range("L122").Value = 20290.5
$ws.Range("M122").Value = -4309
$ws.Range("N122").Value = -25190.5
$ws.Range("H135").Value = 1452.0769
$ws.Range("I135").Value = 175.22223
$ws.Range("J135").Value = 4325
$ws.Range("K135").Value = 1577.00007
$ws.Range("L135").Value = 38925
$ws.Range("M135").Value = 957.9999299999999
$ws.Range("N135").Value = -43995
$ws.Range("H139").Value = 2818.1428
$ws.Range("I139").Value = 2615.6667
$ws.Range("J139").Value = 4033
$ws.Range("K139").Value = 7847.000100000001
$ws.Range("L139").Value = 12099
$ws.Range("M139").Value = -2707.000100000001
$ws.Range("N139").Value = -22379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 12500
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 14000
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = -4707
$ws.Range("N18").Value = -14586
$ws.Range("H35").Value = 9999.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 9999.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 9999.5
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = -10595.5
$ws.Range("H80").Value = 7870.143
$ws.Range("I80").Value = 5028.6665
$ws.Range("J80").Value = 10001.25
$ws.Range("K80").Value = 5028.6665
$ws.Range("L80").Value = 10001.25
$ws.Range("M80").Value = -4030.6665
$ws.Range("N80").Value = -11997.25
$ws.Range("H83").Value = 7870.143
$ws.Range("I83").Value = 5028.6665
$ws.Range("J83").Value = 10001.25
$ws.Range("K83").Value = 25143.3325
$ws.Range("L83").Value = 50006.25
$ws.Range("M83").Value = -20151.3325
$ws.Range("N83").Value = -59990.25
$ws.Range("H122").Value = 2869.1538
$ws.Range("I122").Value = 1811
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 5433
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -2983
$ws.Range("N122").Value = -20650
$ws.Range("H132").Value = 3131.0908
$ws.Range("I132").Value = 2244.2
$ws.Range("K132").Value = 6732.599999999999
$ws.Range("M132").Value = -4202.599999999999
$ws.Range("H136").Value = 42318
$ws.Range("J136").Value = 42318
$ws.Range("L136").Value = 126954
$ws.Range("N136").Value = -132054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 997.1539
$ws.Range("I16").Value = 1137
$ws.Range("K16").Value = 1137
$ws.Range("M16").Value = -967
$ws.Range("H55").Value = 620.125
$ws.Range("I55").Value = 505.26666
$ws.Range("J55").Value = 811.55554
$ws.Range("K55").Value = 505.26666
$ws.Range("L55").Value = 811.55554
$ws.Range("M55").Value = -332.26666
$ws.Range("N55").Value = -1157.55554
$ws.Range("H93").Value = 4059.3157
$ws.Range("I93").Value = 2973.5715
$ws.Range("J93").Value = 7099.4
$ws.Range("K93").Value = 2973.5715
$ws.Range("L93").Value = 7099.4
$ws.Range("M93").Value = -1725.5715
$ws.Range("N93").Value = -9595.4
$ws.Range("H132").Value = 2784.2632
$ws.Range("I132").Value = 2238.2
$ws.Range("K132").Value = 6714.599999999999
$ws.Range("M132").Value = -4184.599999999999
$ws.Range("H136").Value = 5914.0527
$ws.Range("I136").Value = 3297.8
$ws.Range("J136").Value = 15725
$ws.Range("K136").Value = 9893.400000000001
$ws.Range("L136").Value = 47175
$ws.Range("M136").Value = -7343.400000000001
$ws.Range("N136").Value = -52275

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 26498.5
$ws.Range("J45").Value = 26498.5
$ws.Range("L45").Value = 26498.5
$ws.Range("N45").Value = -27480.5
$ws.Range("H81").Value = 4812.3335
$ws.Range("I81").Value = 2903.5833
$ws.Range("J81").Value = 7866.3335
$ws.Range("K81").Value = 5807.1666
$ws.Range("L81").Value = 15732.667
$ws.Range("M81").Value = -4746.1666
$ws.Range("N81").Value = -17854.667
$ws.Range("H84").Value = 4812.3335
$ws.Range("I84").Value = 2903.5833
$ws.Range("J84").Value = 7866.3335
$ws.Range("K84").Value = 29035.833
$ws.Range("L84").Value = 78663.33499999999
$ws.Range("M84").Value = -23731.833
$ws.Range("N84").Value = -89271.33499999999
$ws.Range("H136").Value = 3961.1428
$ws.Range("J136").Value = 6700
$ws.Range("L136").Value = 20100
$ws.Range("N136").Value = -25200
$ws.Range("H137").Value = 86715
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null
